$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.936.95'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.639.46'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.39'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.38'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.17%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.33%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0892'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.873.63'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.638.65'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.30%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.560'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.45%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.56'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.899.72'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '233.10'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0723'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.61'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.31'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.98'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.20%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '150.56'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.94%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.95'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.93%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.25%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.474.23'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.83%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.58%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.89%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.35'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.567'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.927'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +13.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.879'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '69.09'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +6.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.27%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.07%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.23'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.42'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.782.68'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.18%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '87.34'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.64%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0992'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.06%  '
